# Staging.OrganizationType.xlsx: add a "BusinessKey" header column.
#
# Existing header row (row 2) is: A2=Description, B2=OrganizationType_ID.
# The new layout inserts "BusinessKey" as the first header column and
# shifts the rest right by one: A2=BusinessKey, B2=Description,
# C2=OrganizationType_ID.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current header text (Value2 avoids locale/number formatting).
$oldA2 = $ws.Range("A2").Value2
$oldB2 = $ws.Range("B2").Value2

# Shift "OrganizationType_ID" (B2) into the new column C2, copying the
# bold+underline header formatting used by the other header cells.
$ws.Range("C2").Value = $oldB2
$ws.Range("C2").Font.Bold = $true
$ws.Range("C2").Font.Underline = 2

# Shift "Description" (A2) into B2.
$ws.Range("B2").Value = $oldA2

# Insert the new "BusinessKey" header in A2 (keeps A2's existing style).
$ws.Range("A2").Value = "BusinessKey"
